$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (new date, new EBITDA or $null if unchanged)
$updates = @(
    @{ Row = 2;  Date = "2025/11/14"; Ebitda = "4.66" },
    @{ Row = 8;  Date = "2025/11/14"; Ebitda = "7.66" },
    @{ Row = 14; Date = "2025/11/14"; Ebitda = $null },
    @{ Row = 20; Date = "2025/11/14"; Ebitda = "12.13" },
    @{ Row = 26; Date = "2025/11/14"; Ebitda = "9.82" },
    @{ Row = 32; Date = "2025/11/14"; Ebitda = "24.60" },
    @{ Row = 38; Date = "2025/11/14"; Ebitda = $null },
    @{ Row = 44; Date = "2025/11/14"; Ebitda = "10.65" },
    @{ Row = 50; Date = "2025/11/14"; Ebitda = "11.18" },
    @{ Row = 56; Date = "2025/11/14"; Ebitda = "30.71" },
    @{ Row = 62; Date = "2025/11/14"; Ebitda = "11.01" },
    @{ Row = 68; Date = "2025/11/14"; Ebitda = "12.45" },
    @{ Row = 74; Date = "2025/11/14"; Ebitda = "15.14" }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $u.Date
    if ($u.Ebitda -ne $null) {
        $ws.Cells.Item($r, 2).NumberFormat = "@"
        $ws.Cells.Item($r, 2).Value = $u.Ebitda
    }
}
